$d = $word.ActiveDocument
$cmds = $d.list_commands()
Write-Host $cmds
